$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.741.70'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.208.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("E6").Value = '  -4.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.29'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.43%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.402'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.11'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.535.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.23%  '
$ws.Range("E16").Value = '  -2.21%  '
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.187.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.752.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0901'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.67%  '
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.16%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.85%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.76'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.57%  '
$ws.Range("E33").Value = '  -3.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("E35").Value = '  -3.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0647'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.03%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.01%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.35'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.97%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.85%  '
$ws.Range("B40").Value = 'BinanceUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'TerraClassic'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.000238'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0239'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.83%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0954'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.465.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.23%  '
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("E51").Value = '  -5.08%  '
